$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

Set-TextValue "D2" "41.260.09"
Set-TextValue "D3" "2.436.16"
$ws.Range("E3").Value = "  -1.36%  "
Set-TextValue "D5" "318.34"
$ws.Range("E5").Value = "  +0.19%  "
Set-TextValue "D6" "89.85"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -2.73%  "
Set-TextValue "D10" "0.0840"
$ws.Range("E10").Value = "  -2.08%  "
Set-TextValue "D11" "32.17"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("E12").Value = "  -1.67%  "
Set-TextValue "D13" "2.809.01"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("E14").Value = "  -1.64%  "
Set-TextValue "D15" "15.69"
$ws.Range("E15").Value = "  +1.00%  "
Set-TextValue "D16" "2.431.86"
$ws.Range("E16").Value = "  -1.60%  "
Set-TextValue "D17" "0.778"
$ws.Range("E17").Value = "  -1.37%  "
Set-TextValue "D18" "41.162.55"
$ws.Range("E18").Value = "  -0.92%  "
Set-TextValue "D19" "0.0₃0929"
$ws.Range("E19").Value = "  -1.95%  "
Set-TextValue "D20" "6.30"
$ws.Range("E20").Value = "  -1.95%  "
Set-TextValue "D21" "71.98"
$ws.Range("E21").Value = "  +1.25%  "
Set-TextValue "D22" "11.12"
$ws.Range("E22").Value = "  -1.42%  "
Set-TextValue "D23" "236.68"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("E28").Value = "  -2.03%  "
Set-TextValue "D29" "9.64"
$ws.Range("E29").Value = "  -2.04%  "
Set-TextValue "D30" "34.80"
$ws.Range("E30").Value = "  -3.38%  "
Set-TextValue "D31" "156.26"
$ws.Range("E31").Value = "  -3.22%  "
Set-TextValue "D32" "5.29"
$ws.Range("E33").Value = "  +0.07%  "
Set-TextValue "D34" "0.0749"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  -2.87%  "
Set-TextValue "D36" "2.97"
$ws.Range("E36").Value = "  +2.73%  "
Set-TextValue "D37" "16.86"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  -0.57%  "
Set-TextValue "D39" "1.79"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("E40").Value = "  -1.42%  "
Set-TextValue "D41" "3.92"
$ws.Range("E41").Value = "  -0.97%  "
Set-TextValue "D42" "2.003.78"
$ws.Range("E42").Value = "  +0.96%  "
Set-TextValue "D43" "2.22"
$ws.Range("E43").Value = "  -10.21%  "
Set-TextValue "D44" "18.82"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("E45").Value = "  -2.81%  "
Set-TextValue "D46" "2.93"
$ws.Range("E46").Value = "  -1.58%  "
Set-TextValue "D47" "9.56"
$ws.Range("E47").Value = "  +4.31%  "
Set-TextValue "D48" "2.665.48"
$ws.Range("E48").Value = "  -1.63%  "
Set-TextValue "D49" "95.19"
$ws.Range("E49").Value = "  -2.19%  "
Set-TextValue "D50" "73.87"
$ws.Range("E50").Value = "  -0.15%  "
Set-TextValue "D51" "52.14"
$ws.Range("E51").Value = "  -0.26%  "
